$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) updates: force text to avoid numeric coercion of
# values that look like plain decimals, while keeping the default cell
# style (no explicit number format) exactly like the original file.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.301.47'
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.504.77'
$ws.Range('D3').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '590.34'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.41'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.32'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.100.02'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.502.15'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.320.63'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '25.67'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '9.87'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.52'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '393.42'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.572'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.644.20'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '74.64'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.38'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.26'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.48'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.526.07'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '23.49'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.14'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '167.63'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0781'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '25.02'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.344.38'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.893'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '21.15'
$ws.Range('D51').Style = 'Normal'

# Volume(1h) column (E) updates: plain text values, already non-numeric
# (percent sign + padding spaces), so a direct assignment is safe.
$ws.Range('E3').Value = '  -0.54%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('E5').Value = '  +0.92%  '
$ws.Range('E6').Value = '  -0.12%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -0.53%  '
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('E10').Value = '  +2.72%  '
$ws.Range('E11').Value = '  +2.43%  '
$ws.Range('E12').Value = '  -0.56%  '
$ws.Range('E13').Value = '  +1.13%  '
$ws.Range('E14').Value = '  +1.07%  '
$ws.Range('E15').Value = '  -0.56%  '
$ws.Range('E16').Value = '  +0.21%  '
$ws.Range('E17').Value = '  -6.66%  '
$ws.Range('E18').Value = '  +0.82%  '
$ws.Range('E19').Value = '  +2.33%  '
$ws.Range('E20').Value = '  -2.87%  '
$ws.Range('E21').Value = '  +2.72%  '
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('E23').Value = '  -0.59%  '
$ws.Range('E24').Value = '  +0.79%  '
$ws.Range('E25').Value = '  +0.24%  '
$ws.Range('E26').Value = '  -0.20%  '
$ws.Range('E27').Value = '  +0.00%  '
$ws.Range('E28').Value = '  -1.29%  '
$ws.Range('E29').Value = '  +1.38%  '
$ws.Range('E30').Value = '  -2.38%  '
$ws.Range('E31').Value = '  -7.05%  '
$ws.Range('E32').Value = '  -0.30%  '
$ws.Range('E33').Value = '  +5.53%  '
$ws.Range('E35').Value = '  -0.55%  '
$ws.Range('E36').Value = '  -5.28%  '
$ws.Range('E37').Value = '  -0.95%  '
$ws.Range('E38').Value = '  -0.61%  '
$ws.Range('E39').Value = '  +4.34%  '
$ws.Range('E40').Value = '  -0.82%  '
$ws.Range('E41').Value = '  -0.16%  '
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('E43').Value = '  -6.30%  '
$ws.Range('E44').Value = '  -0.22%  '
$ws.Range('E45').Value = '  +3.10%  '
$ws.Range('E46').Value = '  -4.19%  '
$ws.Range('E47').Value = '  -0.79%  '
$ws.Range('E48').Value = '  -5.64%  '
$ws.Range('E49').Value = '  -1.87%  '
$ws.Range('E50').Value = '  -1.47%  '
$ws.Range('E51').Value = '  -1.30%  '
